$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 89.35833500000001
$ws.Range("H2").Value = 268.075005
$ws.Range("I2").Value = 0.9624640326757887
$ws.Range("J2").Value = 0.9624640326757889
$ws.Range("M2").Value = 2.762460333333333
$ws.Range("N2").Value = 8.287381
$ws.Range("O2").Value = 0.2369432190408618
$ws.Range("P2").Value = 0.2369432190408618
$ws.Range("Q2").Value = 246.8488558902117
$ws.Range("R2").Value = 2221.639703011905
$ws.Range("S2").Value = 0.2280493261132506
$ws.Range("T2").Value = 0.2280493261132506
$ws.Range("G3").Value = 89.35833500000001
$ws.Range("H3").Value = 268.075005
$ws.Range("I3").Value = 0.9624640326757887
$ws.Range("J3").Value = 0.9624640326757889
$ws.Range("O3").Value = 0.3694614102513958
$ws.Range("P3").Value = 0.3694614102513958
$ws.Range("Q3").Value = 384.90709624576
$ws.Range("R3").Value = 3464.16386621184
$ws.Range("S3").Value = 0.3555933188286424
$ws.Range("T3").Value = 0.3555933188286425
$ws.Range("G4").Value = 89.35833500000001
$ws.Range("H4").Value = 268.075005
$ws.Range("I4").Value = 0.9624640326757887
$ws.Range("J4").Value = 0.9624640326757889
$ws.Range("M4").Value = 1.776664666666667
$ws.Range("N4").Value = 5.329994000000001
$ws.Range("O4").Value = 0.1523890280691185
$ws.Range("P4").Value = 0.1523890280691185
$ws.Range("Q4").Value = 158.7597964666634
$ws.Range("R4").Value = 1428.83816819997
$ws.Range("S4").Value = 0.1466689584909477
$ws.Range("T4").Value = 0.1466689584909478
$ws.Range("G5").Value = 89.35833500000001
$ws.Range("H5").Value = 268.075005
$ws.Range("I5").Value = 0.9624640326757887
$ws.Range("J5").Value = 0.9624640326757889
$ws.Range("M5").Value = 1.351364
$ws.Range("N5").Value = 4.054092
$ws.Range("O5").Value = 0.1159099127659034
$ws.Range("P5").Value = 0.1159099127659034
$ws.Range("Q5").Value = 120.75563701894
$ws.Range("R5").Value = 1086.80073317046
$ws.Range("S5").Value = 0.1115591220677703
$ws.Range("T5").Value = 0.1115591220677703
$ws.Range("G6").Value = 89.35833500000001
$ws.Range("H6").Value = 268.075005
$ws.Range("I6").Value = 0.9624640326757887
$ws.Range("J6").Value = 0.9624640326757889
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.460799
$ws.Range("N6").Value = 4.382397
$ws.Range("O6").Value = 0.1252964298727204
$ws.Range("P6").Value = 0.1252964298727204
$ws.Range("Q6").Value = 130.534566409665
$ws.Range("R6").Value = 1174.811097686985
$ws.Range("S6").Value = 0.1205933071751776
$ws.Range("T6").Value = 0.1205933071751777
$ws.Range("I7").Value = 0.001854741667334279
$ws.Range("J7").Value = 0.001854741667334279
$ws.Range("M7").Value = 2.762460333333333
$ws.Range("N7").Value = 8.287381
$ws.Range("O7").Value = 0.2369432190408618
$ws.Range("P7").Value = 0.2369432190408618
$ws.Range("Q7").Value = 0.4756965902201112
$ws.Range("R7").Value = 4.281269311981001
$ws.Range("S7").Value = 0.0004394684611473993
$ws.Range("T7").Value = 0.0004394684611473993
$ws.Range("I8").Value = 0.001854741667334279
$ws.Range("J8").Value = 0.001854741667334279
$ws.Range("O8").Value = 0.3694614102513958
$ws.Range("P8").Value = 0.3694614102513958
$ws.Range("S8").Value = 0.0006852554720653479
$ws.Range("T8").Value = 0.000685255472065348
$ws.Range("I9").Value = 0.001854741667334279
$ws.Range("J9").Value = 0.001854741667334279
$ws.Range("M9").Value = 1.776664666666667
$ws.Range("N9").Value = 5.329994000000001
$ws.Range("O9").Value = 0.1523890280691185
$ws.Range("P9").Value = 0.1523890280691185
$ws.Range("Q9").Value = 0.3059422478215557
$ws.Range("R9").Value = 2.753480230394001
$ws.Range("S9").Value = 0.0002826422800043671
$ws.Range("T9").Value = 0.0002826422800043671
$ws.Range("I10").Value = 0.001854741667334279
$ws.Range("J10").Value = 0.001854741667334279
$ws.Range("M10").Value = 1.351364
$ws.Range("N10").Value = 4.054092
$ws.Range("O10").Value = 0.1159099127659034
$ws.Range("P10").Value = 0.1159099127659034
$ws.Range("Q10").Value = 0.2327053312546667
$ws.Range("R10").Value = 2.094347981292
$ws.Range("S10").Value = 0.0002149829448640026
$ws.Range("T10").Value = 0.0002149829448640026
$ws.Range("I11").Value = 0.001854741667334279
$ws.Range("J11").Value = 0.001854741667334279
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.460799
$ws.Range("N11").Value = 4.382397
$ws.Range("O11").Value = 0.1252964298727204
$ws.Range("P11").Value = 0.1252964298727204
$ws.Range("Q11").Value = 0.2515500747330001
$ws.Range("R11").Value = 2.263950672597
$ws.Range("S11").Value = 0.000232392509253162
$ws.Range("T11").Value = 0.000232392509253162
$ws.Range("G12").Value = 1.963978
$ws.Range("H12").Value = 5.891934
$ws.Range("I12").Value = 0.02115368628977398
$ws.Range("J12").Value = 0.02115368628977398
$ws.Range("M12").Value = 2.762460333333333
$ws.Range("N12").Value = 8.287381
$ws.Range("O12").Value = 0.2369432190408618
$ws.Range("P12").Value = 0.2369432190408618
$ws.Range("Q12").Value = 5.425411320539333
$ws.Range("R12").Value = 48.828701884854
$ws.Range("S12").Value = 0.005012222524079591
$ws.Range("T12").Value = 0.005012222524079591
$ws.Range("G13").Value = 1.963978
$ws.Range("H13").Value = 5.891934
$ws.Range("I13").Value = 0.02115368628977398
$ws.Range("J13").Value = 0.02115368628977398
$ws.Range("O13").Value = 0.3694614102513958
$ws.Range("P13").Value = 0.3694614102513958
$ws.Range("Q13").Value = 8.459748819967999
$ws.Range("R13").Value = 76.13773937971199
$ws.Range("S13").Value = 0.007815470768635509
$ws.Range("T13").Value = 0.007815470768635511
$ws.Range("G14").Value = 1.963978
$ws.Range("H14").Value = 5.891934
$ws.Range("I14").Value = 0.02115368628977398
$ws.Range("J14").Value = 0.02115368628977398
$ws.Range("M14").Value = 1.776664666666667
$ws.Range("N14").Value = 5.329994000000001
$ws.Range("O14").Value = 0.1523890280691185
$ws.Range("P14").Value = 0.1523890280691185
$ws.Range("Q14").Value = 3.489330318710667
$ws.Range("R14").Value = 31.40397286839601
$ws.Range("S14").Value = 0.003223589693777694
$ws.Range("T14").Value = 0.003223589693777694
$ws.Range("G15").Value = 1.963978
$ws.Range("H15").Value = 5.891934
$ws.Range("I15").Value = 0.02115368628977398
$ws.Range("J15").Value = 0.02115368628977398
$ws.Range("M15").Value = 1.351364
$ws.Range("N15").Value = 4.054092
$ws.Range("O15").Value = 0.1159099127659034
$ws.Range("P15").Value = 0.1159099127659034
$ws.Range("Q15").Value = 2.654049165992
$ws.Range("R15").Value = 23.886442493928
$ws.Range("S15").Value = 0.002451921932524989
$ws.Range("T15").Value = 0.002451921932524989
$ws.Range("G16").Value = 1.963978
$ws.Range("H16").Value = 5.891934
$ws.Range("I16").Value = 0.02115368628977398
$ws.Range("J16").Value = 0.02115368628977398
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.460799
$ws.Range("N16").Value = 4.382397
$ws.Range("O16").Value = 0.1252964298727204
$ws.Range("P16").Value = 0.1252964298727204
$ws.Range("Q16").Value = 2.868977098422
$ws.Range("R16").Value = 25.820793885798
$ws.Range("S16").Value = 0.002650481370756192
$ws.Range("T16").Value = 0.002650481370756193
$ws.Range("G17").Value = 0.16825
$ws.Range("H17").Value = 0.50475
$ws.Range("I17").Value = 0.001812193272151965
$ws.Range("J17").Value = 0.001812193272151965
$ws.Range("M17").Value = 2.762460333333333
$ws.Range("N17").Value = 8.287381
$ws.Range("O17").Value = 0.2369432190408618
$ws.Range("P17").Value = 0.2369432190408618
$ws.Range("Q17").Value = 0.4647839510833334
$ws.Range("R17").Value = 4.183055559750001
$ws.Range("S17").Value = 0.0004293869074278791
$ws.Range("T17").Value = 0.0004293869074278791
$ws.Range("G18").Value = 0.16825
$ws.Range("H18").Value = 0.50475
$ws.Range("I18").Value = 0.001812193272151965
$ws.Range("J18").Value = 0.001812193272151965
$ws.Range("O18").Value = 0.3694614102513958
$ws.Range("P18").Value = 0.3694614102513958
$ws.Range("Q18").Value = 0.7247294719999999
$ws.Range("R18").Value = 6.522565247999999
$ws.Range("S18").Value = 0.0006695354819773564
$ws.Range("T18").Value = 0.0006695354819773565
$ws.Range("G19").Value = 0.16825
$ws.Range("H19").Value = 0.50475
$ws.Range("I19").Value = 0.001812193272151965
$ws.Range("J19").Value = 0.001812193272151965
$ws.Range("M19").Value = 1.776664666666667
$ws.Range("N19").Value = 5.329994000000001
$ws.Range("O19").Value = 0.1523890280691185
$ws.Range("P19").Value = 0.1523890280691185
$ws.Range("Q19").Value = 0.2989238301666667
$ws.Range("R19").Value = 2.690314471500001
$ws.Range("S19").Value = 0.0002761583714166335
$ws.Range("T19").Value = 0.0002761583714166335
$ws.Range("G20").Value = 0.16825
$ws.Range("H20").Value = 0.50475
$ws.Range("I20").Value = 0.001812193272151965
$ws.Range("J20").Value = 0.001812193272151965
$ws.Range("M20").Value = 1.351364
$ws.Range("N20").Value = 4.054092
$ws.Range("O20").Value = 0.1159099127659034
$ws.Range("P20").Value = 0.1159099127659034
$ws.Range("Q20").Value = 0.227366993
$ws.Range("R20").Value = 2.046302937
$ws.Range("S20").Value = 0.0002100511640900913
$ws.Range("T20").Value = 0.0002100511640900913
$ws.Range("G21").Value = 0.16825
$ws.Range("H21").Value = 0.50475
$ws.Range("I21").Value = 0.001812193272151965
$ws.Range("J21").Value = 0.001812193272151965
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 1.460799
$ws.Range("N21").Value = 4.382397
$ws.Range("O21").Value = 0.1252964298727204
$ws.Range("P21").Value = 0.1252964298727204
$ws.Range("Q21").Value = 0.24577943175
$ws.Range("R21").Value = 2.21201488575
$ws.Range("S21").Value = 0.0002270613472400044
$ws.Range("T21").Value = 0.0002270613472400044
$ws.Range("G22").Value = 1.180534666666667
$ws.Range("H22").Value = 3.541604
$ws.Range("I22").Value = 0.01271534609495094
$ws.Range("J22").Value = 0.01271534609495094
$ws.Range("M22").Value = 2.762460333333333
$ws.Range("N22").Value = 8.287381
$ws.Range("O22").Value = 0.2369432190408618
$ws.Range("P22").Value = 0.2369432190408618
$ws.Range("Q22").Value = 3.261180188791556
$ws.Range("R22").Value = 29.350621699124
$ws.Range("S22").Value = 0.003012815034956327
$ws.Range("T22").Value = 0.003012815034956327
$ws.Range("G23").Value = 1.180534666666667
$ws.Range("H23").Value = 3.541604
$ws.Range("I23").Value = 0.01271534609495094
$ws.Range("J23").Value = 0.01271534609495094
$ws.Range("O23").Value = 0.3694614102513958
$ws.Range("P23").Value = 0.3694614102513958
$ws.Range("Q23").Value = 5.085101133141332
$ws.Range("R23").Value = 45.765910198272
$ws.Range("S23").Value = 0.004697829700075153
$ws.Range("T23").Value = 0.004697829700075153
$ws.Range("G24").Value = 1.180534666666667
$ws.Range("H24").Value = 3.541604
$ws.Range("I24").Value = 0.01271534609495094
$ws.Range("J24").Value = 0.01271534609495094
$ws.Range("M24").Value = 1.776664666666667
$ws.Range("N24").Value = 5.329994000000001
$ws.Range("O24").Value = 0.1523890280691185
$ws.Range("P24").Value = 0.1523890280691185
$ws.Range("Q24").Value = 2.097414230041778
$ws.Range("R24").Value = 18.876728070376
$ws.Range("S24").Value = 0.001937679232972035
$ws.Range("T24").Value = 0.001937679232972035
$ws.Range("G25").Value = 1.180534666666667
$ws.Range("H25").Value = 3.541604
$ws.Range("I25").Value = 0.01271534609495094
$ws.Range("J25").Value = 0.01271534609495094
$ws.Range("M25").Value = 1.351364
$ws.Range("N25").Value = 4.054092
$ws.Range("O25").Value = 0.1159099127659034
$ws.Range("P25").Value = 0.1159099127659034
$ws.Range("Q25").Value = 1.595332049285333
$ws.Range("R25").Value = 14.357988443568
$ws.Range("S25").Value = 0.001473834656654034
$ws.Range("T25").Value = 0.001473834656654034
$ws.Range("G26").Value = 1.180534666666667
$ws.Range("H26").Value = 3.541604
$ws.Range("I26").Value = 0.01271534609495094
$ws.Range("J26").Value = 0.01271534609495094
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 1.460799
$ws.Range("N26").Value = 4.382397
$ws.Range("O26").Value = 0.1252964298727204
$ws.Range("P26").Value = 0.1252964298727204
$ws.Range("Q26").Value = 1.724523860532
$ws.Range("R26").Value = 15.520714744788
$ws.Range("S26").Value = 0.00159318747029339
$ws.Range("T26").Value = 0.00159318747029339
